$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 92: blank separator row between Feb 07 block and the new Feb 10 block.
# Reuse the exact same "day separator" formatting already used elsewhere in the
# sheet (e.g. row 29) so no new style entries are created.
$ws.Range("A29:C29").Copy()
$ws.Range("A92:C92").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 94 (B column) needs the wrapped / multi-line style already used for
# long task descriptions elsewhere (e.g. B4). Copy that formatting over before
# writing the value.
$ws.Range("B4").Copy()
$ws.Range("B94").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Fill in the new Feb 10 2020 timesheet entries (rows 93-100).
# Values are written in this particular order so that the generated shared
# string table lines up with how the workbook was actually authored.
$ws.Range("A93").Value = "Feb 10 10:00 to 11:00"
$ws.Range("A94").Value = "Feb 10 11:00 to 12:00"
$ws.Range("B94").Value = "Instaling project setup on my local machine but getting problems on installing`ndependencies."
$ws.Range("A95").Value = "Feb 10 12:00 to 13:00"
$ws.Range("B95").Value = "Discussed with sujata mam."
$ws.Range("B93").Value = "Watched video that was recommeded by client."
$ws.Range("A96").Value = "Feb 10 13:00 to 14:00"
$ws.Range("B96").Value = "Converted mp3 file to wav file"
$ws.Range("A97").Value = "Feb 10 14:00 to 15:00"
$ws.Range("B98").Value = "Program giving error while doing training and testing split. 'numpy memory error'"
$ws.Range("A98").Value = "Feb 10 15:00 to 16:00"
$ws.Range("A99").Value = "Feb 10 16:00 to 17:00"
$ws.Range("B99").Value = "Looked at all the videos and understand the procedure of audio seperation."
$ws.Range("A100").Value = "Feb 10 17:00 to 18:00"
$ws.Range("B100").Value = "Working on making 24 buckets rather than 1 day bin bucket."

# Column C ("Location") stays "Infimetrics" for every entry of this day, and
# B97 is another lunch break, reusing the existing shared strings.
$ws.Range("C93").Value = "Infimetrics"
$ws.Range("C94").Value = "Infimetrics"
$ws.Range("C95").Value = "Infimetrics"
$ws.Range("C96").Value = "Infimetrics"
$ws.Range("B97").Value = "Lunch"
$ws.Range("C97").Value = "Infimetrics"
$ws.Range("C98").Value = "Infimetrics"
$ws.Range("C99").Value = "Infimetrics"
$ws.Range("C100").Value = "Infimetrics"

# Row 94 wraps onto two lines in the real workbook (ht="30").
$ws.Rows.Item(94).RowHeight = 30

# Match the scrolled/selected view state left behind after the edit.
$ws.Range("C101").Select()
$excel.ActiveWindow.ScrollRow = 88
$excel.ActiveWindow.ScrollColumn = 1
